{"js": "// Office.js (Word JavaScript API) script.\n// Applies the \"Support for starred items\" edit described by the diff:\n//  1. Splits \"reading\" in the \"Easy navigation...\" bullet and inserts a\n//     collapsed \"_GoBack\" bookmark between \"read\" and \"ing\" (bookmark\n//     moved here from the paragraph that gets removed below).\n//  2. \"New in 4.4\" -> \"New in 5.0\"\n//  3. The language bullet becomes \"New Languages - French, Bahasa,\n//     Brazilian Portuguese and Swedish\" (several runs, matching the\n//     original run-splitting around \"Bahasa\").\n//  4. The old \"Experimental languages - ... Swedish and Turkish\" bullet\n//     is removed (its content + _GoBack bookmark got folded away).\n//  5. A new blank paragraph is left where that bullet used to be.\n//  6. \"Fixed a crash for large local cache\" gets a trailing \".\" run.\n//  7. \"Better progress status\" gets a trailing \".\" run.\n//  8. New bullet \"Ability to check for updates.\" is added right after.\n\nconst OOXML_NS = 'xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"';\n\nfunction wrapParagraphOoxml(innerParagraphXml) {\n  return (\n    '<pkg:package ' + OOXML_NS + '>' +\n      '<pkg:part pkg:name=\"/word/document.xml\" ' +\n        'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n          '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n            '<w:body>' + innerParagraphXml + '</w:body>' +\n          '</w:document>' +\n        '</pkg:xmlData>' +\n      '</pkg:part>' +\n    '</pkg:package>'\n  );\n}\n\nconst body = context.document.body;\n\n// --- 1. Move the _GoBack bookmark into the \"Easy navigation\" bullet,\n//        splitting \"reading\" into \"read\" + \"ing\". ----------------------\n// Search precisely for the \"read\" substring that immediately precedes\n// \"ing a feed\" and bookmark its (collapsed) end point.\nconst readResults = body.search(\"while read\", { matchCase: true });\nawait context.sync();\nif (readResults.items.length !== 1) {\n  throw new Error(\"Expected exactly one 'while read' match, got \" + readResults.items.length);\n}\nconst readEnd = readResults.items[0].getRange(\"End\");\nreadEnd.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 2. \"New in 4.4\" -> \"New in 5.0\" ------------------------------------\nconst verResults = body.search(\"New in 4.4\", { matchCase: true });\nawait context.sync();\nif (verResults.items.length !== 1) {\n  throw new Error(\"Expected exactly one 'New in 4.4' match, got \" + verResults.items.length);\n}\nverResults.items[0].insertText(\"New in 5.0\", Word.InsertLocation.replace);\nawait context.sync();\n\n// --- 3 & 4. Rewrite the languages bullet and drop the old one -----------\n// Replace \" and Brazilian Portuguese\" (note leading space) with\n// \", \" + \"Brazilian Portuguese\" + \" and Swedish\" as three distinct runs,\n// matching the original's \"Bahasa\" run-splitting style.\nconst langResults = body.search(\" and Brazilian Portuguese\", { matchCase: true });\nawait context.sync();\nif (langResults.items.length !== 1) {\n  throw new Error(\"Expected exactly one ' and Brazilian Portuguese' match, got \" + langResults.items.length);\n}\nlangResults.items[0].insertOoxml(\n  wrapParagraphOoxml(\n    '<w:p>' +\n      '<w:r><w:t xml:space=\"preserve\">, </w:t></w:r>' +\n      '<w:r><w:t>Brazilian Portuguese</w:t></w:r>' +\n      '<w:r><w:t xml:space=\"preserve\"> and Swedish</w:t></w:r>' +\n    '</w:p>'\n  ),\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// Remove the now-obsolete \"Experimental languages ... Swedish and Turkish\"\n// paragraph (its bookmark already migrated to the bullet above), but keep\n// an empty paragraph in its place (matches the diff: a blank <w:p/> stays\n// between the languages bullet and \"Fixed a crash...\").\nbody.paragraphs.load(\"items/text\");\nawait context.sync();\nlet experimentalPara = null;\nfor (const p of body.paragraphs.items) {\n  if (p.text.indexOf(\"Experimental languages\") === 0) {\n    experimentalPara = p;\n    break;\n  }\n}\nif (!experimentalPara) {\n  throw new Error(\"Could not find the 'Experimental languages' paragraph\");\n}\n// The pre-existing blank paragraph that originally followed the\n// languages bullets (between it and \"Fixed a crash...\") is untouched,\n// so simply deleting this paragraph leaves exactly one blank line.\nexperimentalPara.delete();\nawait context.sync();\n\n// --- 6. \"Fixed a crash for large local cache\" gets a trailing \".\" run ---\nconst cacheResults = body.search(\"Fixed a crash for large local cache\", { matchCase: true });\nawait context.sync();\nif (cacheResults.items.length !== 1) {\n  throw new Error(\"Expected exactly one 'Fixed a crash...' match, got \" + cacheResults.items.length);\n}\ncacheResults.items[0].getRange(\"End\").insertOoxml(\n  wrapParagraphOoxml('<w:p><w:r><w:t>.</w:t></w:r></w:p>'),\n  Word.InsertLocation.end\n);\nawait context.sync();\n\n// --- 7 & 8. \"Better progress status\" gets a trailing \".\" run, and a new\n//            \"Ability to check for updates.\" bullet follows it. ---------\nconst progressResults = body.search(\"Better progress status\", { matchCase: true });\nawait context.sync();\nif (progressResults.items.length !== 1) {\n  throw new Error(\"Expected exactly one 'Better progress status' match, got \" + progressResults.items.length);\n}\nprogressResults.items[0].getRange(\"End\").insertOoxml(\n  wrapParagraphOoxml('<w:p><w:r><w:t>.</w:t></w:r></w:p>'),\n  Word.InsertLocation.end\n);\nawait context.sync();\n\nbody.paragraphs.load(\"items/text\");\nawait context.sync();\nlet progressPara = null;\nfor (const p of body.paragraphs.items) {\n  if (p.text.trim() === \"Better progress status.\") {\n    progressPara = p;\n    break;\n  }\n}\nif (!progressPara) {\n  throw new Error(\"Could not find the 'Better progress status.' paragraph\");\n}\nprogressPara.insertParagraph(\"Ability to check for updates.\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the \"Support for starred items\" edit described by the diff:\n#  1. Splits \"reading\" in the \"Easy navigation...\" bullet and inserts a\n#     collapsed \"_GoBack\" bookmark between \"read\" and \"ing\" (bookmark\n#     moved here from the paragraph that gets removed below).\n#  2. \"New in 4.4\" -> \"New in 5.0\"\n#  3. The language bullet becomes \"New Languages - French, Bahasa,\n#     Brazilian Portuguese and Swedish\" (several runs, matching the\n#     original run-splitting around \"Bahasa\").\n#  4. The old \"Experimental languages - ... Swedish and Turkish\" bullet\n#     is removed (its content + _GoBack bookmark got folded away).\n#  5. A blank paragraph remains where that bullet used to be.\n#  6. \"Fixed a crash for large local cache\" gets a trailing \".\" run.\n#  7. \"Better progress status\" gets a trailing \".\" run.\n#  8. New bullet \"Ability to check for updates.\" is added right after.\n\n$d = $word.ActiveDocument\n\nfunction New-PkgXml([string]$innerParagraphXml) {\n    return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" ' +\n        'pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' + $innerParagraphXml + '</w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n}\n\nfunction Get-ParagraphByPrefix([string]$prefix) {\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text.StartsWith($prefix)) {\n            return $p\n        }\n    }\n    return $null\n}\n\n# --- 1. Move the _GoBack bookmark into the \"Easy navigation\" bullet, -----\n#        splitting \"reading\" into \"read\" + \"ing\".\n$rng = $d.Content\n$found = $rng.Find.Execute(\"while read\")\nif (-not $found) {\n    throw \"Could not find 'while read'\"\n}\n$bookmarkPoint = $d.Range($rng.End, $rng.End)\n$null = $d.Bookmarks.Add(\"_GoBack\", $bookmarkPoint)\n\n# --- 2. \"New in 4.4\" -> \"New in 5.0\" --------------------------------------\n$verPara = Get-ParagraphByPrefix(\"New in 4.4\")\nif (-not $verPara) {\n    throw \"Could not find the 'New in 4.4' paragraph\"\n}\n$verRng = $verPara.Range\n$verRng.Text = \"\"\n$null = $verRng.InsertXML((New-PkgXml '<w:p><w:r><w:t>New in 5.0</w:t></w:r></w:p>'))\n\n# --- 3. Rewrite the languages bullet --------------------------------------\n$langPara = Get-ParagraphByPrefix(\"New Languages\")\nif (-not $langPara) {\n    throw \"Could not find the 'New Languages' paragraph\"\n}\n$langRng = $langPara.Range\n$langRng.Text = \"\"\n$null = $langRng.InsertXML((New-PkgXml (\n    '<w:p>' +\n    '<w:r><w:t xml:space=\"preserve\">New Languages - French, </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Bahasa</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\">, </w:t></w:r>' +\n    '<w:r><w:t>Brazilian Portuguese</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> and Swedish</w:t></w:r>' +\n    '</w:p>'\n)))\n\n# --- 4 & 5. Remove the now-obsolete \"Experimental languages...\" bullet ---\n#            (a blank paragraph already follows it, so deleting this one\n#             paragraph leaves exactly one blank line behind).\n$expPara = Get-ParagraphByPrefix(\"Experimental languages\")\nif (-not $expPara) {\n    throw \"Could not find the 'Experimental languages' paragraph\"\n}\n$null = $expPara.Range.Delete()\n\n# --- 6. \"Fixed a crash for large local cache\" gets a trailing \".\" run ----\n$cachePara = Get-ParagraphByPrefix(\"Fixed a crash for large local cache\")\nif (-not $cachePara) {\n    throw \"Could not find the 'Fixed a crash...' paragraph\"\n}\n$cacheRng = $cachePara.Range\n$cacheRng.Text = \"\"\n$null = $cacheRng.InsertXML((New-PkgXml '<w:p><w:r><w:t>Fixed a crash for large local cache</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>'))\n\n# --- 7. \"Better progress status\" gets a trailing \".\" run -----------------\n$progressPara = Get-ParagraphByPrefix(\"Better progress status\")\nif (-not $progressPara) {\n    throw \"Could not find the 'Better progress status' paragraph\"\n}\n$progressRng = $progressPara.Range\n$progressRng.Text = \"\"\n$null = $progressRng.InsertXML((New-PkgXml '<w:p><w:r><w:t>Better progress status</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p>'))\n\n# --- 8. New bullet \"Ability to check for updates.\" follows it ------------\n$progressPara = Get-ParagraphByPrefix(\"Better progress status\")\nif (-not $progressPara) {\n    throw \"Could not find the 'Better progress status' paragraph (post-edit)\"\n}\n$null = $progressPara.Range.InsertParagraphAfter()\n$newPara = $progressPara.Next()\n$newPara.Range.Text = \"Ability to check for updates.\"\n"}
